# US_93: Quizzes storage to database
$wb = $excel.ActiveWorkbook

# Rename the "Expression" sheet to "Quizzes"
$ws = $wb.Worksheets.Item("Expression")
$ws.Name = "Quizzes"

# Replace the header row
$ws.Range("A1").Value = "Grade"
$ws.Range("B1").Value = "Question"
$ws.Range("C1").Value = "Answer"

# Fill in quiz data grouped by grade
$data = @(
    @("1", "3+2", "5"),
    @("1", "6+4", "10"),
    @("1", "8+3", "11"),
    @("4", "3+2", "5"),
    @("4", "6+4", "10"),
    @("4", "8+3", "11")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = "'" + $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = "'" + $data[$i][2]
}
